$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataText")

# --- Rerun of the APS (Annual Population Survey / nomis) extract ---
# The employment-rate / employment-volume rows (2-9) carried the caveat
# text for the previous data window ("Jul 2023-Jun 2024..."). The refreshed
# nomis extract covers "Oct 2023-Sep 2024...", so update the LatestPeriod
# (column B) caveat text for each of those rows to match.
$newPeriodText = "Oct 2023-Sep 2024.  North East and North of Tyne MCA have temporarily been removed, on account of North of Tyne being integrated into North East from now on"

foreach ($r in 2..9) {
    $ws.Cells.Item($r, 2).Value = $newPeriodText
}

# --- Update the view state left after the edit (active cell/selection) ---
$ws.Range("E8").Select()
